$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('face/face021.jpg', 'rechnen', 'face'),
    @('face/face020.jpg', 'stoppen', 'face'),
    @('house/house017.jpg', 'kennen', 'house'),
    @('face/face009.jpg', 'quellen', 'face'),
    @('house/house019.jpg', 'prüfen', 'house'),
    @('house/house005.jpg', 'achten', 'house'),
    @('house/house020.jpg', 'piepen', 'house'),
    @('house/house016.jpg', 'leeren', 'house'),
    @('house/house004.jpg', 'mögen', 'house'),
    @('face/face022.jpg', 'duschen', 'face'),
    @('face/face015.jpg', 'wählen', 'face'),
    @('house/house028.jpg', 'hören', 'house'),
    @('face/face028.jpg', 'süßen', 'face'),
    @('house/house008.jpg', 'fließen', 'house'),
    @('face/face000.jpg', 'bergen', 'face'),
    @('face/face025.jpg', 'reisen', 'face'),
    @('face/face029.jpg', 'zögern', 'face'),
    @('face/face024.jpg', 'binden', 'face'),
    @('face/face014.jpg', 'heben', 'face'),
    @('house/house012.jpg', 'betteln', 'house'),
    @('house/house030.jpg', 'ärgern', 'house'),
    @('house/house023.jpg', 'wachsen', 'house'),
    @('house/house006.jpg', 'tollen', 'house'),
    @('face/face012.jpg', 'spenden', 'face'),
    @('house/house018.jpg', 'herrschen', 'house'),
    @('house/house009.jpg', 'nullen', 'house'),
    @('face/face011.jpg', 'öffnen', 'face'),
    @('house/house014.jpg', 'hassen', 'house'),
    @('face/face005.jpg', 'bauen', 'face'),
    @('face/face013.jpg', 'holen', 'face'),
    @('face/face023.jpg', 'zielen', 'face'),
    @('house/house002.jpg', 'wehen', 'house'),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value = $data[$i][0]
    $ws.Cells.Item($r, 3).Value = $data[$i][1]
    $ws.Cells.Item($r, 4).Value = $data[$i][2]
}